$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty row 8 with a new page link entry (no shifting of other rows)
$ws.Range("A8").Value = "ModeratorPanelOwnersPage"
$ws.Range("B8").Value = "/moderator/owners"
$ws.Range("C8").Value = "Anton Tsvihun"

$ws.Range("C8").Select()
